$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 328
$ws.Range("I28").Value = 358.14285
$ws.Range("J28").Value = 257.66666
$ws.Range("K28").Value = 358.14285
$ws.Range("L28").Value = 257.66666
$ws.Range("M28").Value = 126.85715
$ws.Range("N28").Value = -1227.66666
$ws.Range("H32").Value = 4998
$ws.Range("I32").Value = 4998
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 4998
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -4672
$ws.Range("N32").ClearContents()
$ws.Range("H40").Value = 2222.111
$ws.Range("I40").Value = 2000
$ws.Range("J40").Value = 2499.75
$ws.Range("K40").Value = 2000
$ws.Range("L40").Value = 2499.75
$ws.Range("M40").Value = -1825
$ws.Range("N40").Value = -2849.75
$ws.Range("H70").Value = 10708.363
$ws.Range("I70").Value = 710
$ws.Range("J70").Value = 16421.715
$ws.Range("K70").Value = 2130
$ws.Range("L70").Value = 49265.145
$ws.Range("M70").Value = -1860
$ws.Range("N70").Value = -49805.145
$ws.Range("H73").Value = 10708.363
$ws.Range("I73").Value = 710
$ws.Range("J73").Value = 16421.715
$ws.Range("K73").Value = 2130
$ws.Range("L73").Value = 49265.145
$ws.Range("M73").Value = -1194
$ws.Range("N73").Value = -51137.145
$ws.Range("H80").Value = 514.46155
$ws.Range("I80").Value = 355.57144
$ws.Range("J80").Value = 699.8333
$ws.Range("K80").Value = 1066.71432
$ws.Range("L80").Value = 2099.4999
$ws.Range("M80").Value = -68.71432000000004
$ws.Range("H83").Value = 514.46155
$ws.Range("I83").Value = 355.57144
$ws.Range("J83").Value = 699.8333
$ws.Range("K83").Value = 3200.14296
$ws.Range("L83").Value = 6298.4997
$ws.Range("M83").Value = 1791.85704
$ws.Range("H92").Value = 818.26666
$ws.Range("I92").Value = 898.4
$ws.Range("J92").Value = 658
$ws.Range("K92").Value = 898.4
$ws.Range("L92").Value = 658
$ws.Range("M92").Value = 349.6
$ws.Range("N92").Value = -3154
$ws.Range("H106").Value = 21496
$ws.Range("I106").Value = 26550.223
$ws.Range("J106").Value = 6333.3335
$ws.Range("K106").Value = 26550.223
$ws.Range("L106").Value = 6333.3335
$ws.Range("M106").Value = -25919.223
$ws.Range("H115").Value = 2551
$ws.Range("I115").Value = 2551
$ws.Range("J115").Value = 0
$ws.Range("K115").Value = 7653
$ws.Range("L115").Value = 0
$ws.Range("M115").Value = -6086
$ws.Range("H118").Value = 0
$ws.Range("I118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("K118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").ClearContents()
$ws.Range("H137").Value = 3243.2666
$ws.Range("I137").Value = 1842.3334
$ws.Range("J137").Value = 4177.222
$ws.Range("K137").Value = 5527.0002
$ws.Range("L137").Value = 12531.666
$ws.Range("M137").Value = -2977.0002
$ws.Range("N137").Value = -17631.666

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H12").Value = 5048.3335
$ws.Range("I12").Value = 15000
$ws.Range("J12").Value = 72.5
$ws.Range("K12").Value = 15000
$ws.Range("L12").Value = 72.5
$ws.Range("M12").Value = -14827
$ws.Range("N12").Value = -418.5
$ws.Range("H61").Value = 1405.7742
$ws.Range("I61").Value = 1402.7333
$ws.Range("J61").Value = 1497
$ws.Range("K61").Value = 1402.7333
$ws.Range("L61").Value = 1497
$ws.Range("M61").Value = -1190.7333
$ws.Range("N61").Value = -1921
$ws.Range("H63").Value = 7317.727
$ws.Range("I63").Value = 6373.75
$ws.Range("J63").Value = 7857.143
$ws.Range("K63").Value = 6373.75
$ws.Range("L63").Value = 7857.143
$ws.Range("M63").Value = -5687.75
$ws.Range("N63").Value = -9229.143
$ws.Range("H66").Value = 7317.727
$ws.Range("I66").Value = 6373.75
$ws.Range("J66").Value = 7857.143
$ws.Range("K66").Value = 31868.75
$ws.Range("L66").Value = 39285.715
$ws.Range("M66").Value = -28436.75
$ws.Range("N66").Value = -46149.715
$ws.Range("H97").Value = 751.125
$ws.Range("I97").Value = 722.7143
$ws.Range("J97").Value = 950
$ws.Range("K97").Value = 722.7143
$ws.Range("L97").Value = 950
$ws.Range("M97").Value = -226.7143
$ws.Range("H122").Value = 387291.84
$ws.Range("I122").Value = 528368.0600000001
$ws.Range("J122").Value = 4370.7144
$ws.Range("K122").Value = 1585104.18
$ws.Range("L122").Value = 13112.1432
$ws.Range("M122").Value = -1582654.18
$ws.Range("H136").Value = 1405.7742
$ws.Range("I136").Value = 1402.7333
$ws.Range("J136").Value = 1497
$ws.Range("K136").Value = 4208.199900000001
$ws.Range("L136").Value = 4491
$ws.Range("M136").Value = -1658.199900000001
$ws.Range("N136").Value = -9591

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 536.75
$ws.Range("I80").Value = 726
$ws.Range("J80").Value = 347.5
$ws.Range("K80").Value = 726
$ws.Range("L80").Value = 347.5
$ws.Range("M80").Value = 272
$ws.Range("H83").Value = 536.75
$ws.Range("I83").Value = 726
$ws.Range("J83").Value = 347.5
$ws.Range("K83").Value = 3630
$ws.Range("L83").Value = 1737.5
$ws.Range("M83").Value = 1362
$ws.Range("H105").Value = 3990.5173
$ws.Range("I105").Value = 3401.5881
$ws.Range("J105").Value = 4824.8335
$ws.Range("K105").Value = 3401.5881
$ws.Range("L105").Value = 4824.8335
$ws.Range("M105").Value = -1654.5881
$ws.Range("N105").Value = -8318.833500000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 556.1111
$ws.Range("I16").Value = 1005
$ws.Range("J16").Value = 500
$ws.Range("K16").Value = 1005
$ws.Range("L16").Value = 500
$ws.Range("M16").Value = -718
$ws.Range("H22").Value = 348
$ws.Range("I22").Value = 350.16666
$ws.Range("J22").Value = 343.66666
$ws.Range("K22").Value = 350.16666
$ws.Range("L22").Value = 343.66666
$ws.Range("M22").Value = -0.1666599999999789
$ws.Range("N22").Value = -1043.66666
$ws.Range("H44").Value = 12016
$ws.Range("I44").Value = 7532
$ws.Range("J44").Value = 16500
$ws.Range("K44").Value = 7532
$ws.Range("L44").Value = 16500
$ws.Range("M44").Value = -7090
$ws.Range("H52").Value = 98999
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 98999
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 98999
$ws.Range("N52").Value = -99587
$ws.Range("H58").Value = 3469.2942
$ws.Range("I58").Value = 1935.1111
$ws.Range("J58").Value = 5195.25
$ws.Range("K58").Value = 1935.1111
$ws.Range("L58").Value = 5195.25
$ws.Range("M58").Value = -1732.1111
$ws.Range("H113").Value = 556.1111
$ws.Range("I113").Value = 1005
$ws.Range("J113").Value = 500
$ws.Range("K113").Value = 1005
$ws.Range("L113").Value = 500
$ws.Range("M113").Value = 1165
$ws.Range("H132").Value = 2442.0967
$ws.Range("I132").Value = 1993.75
$ws.Range("J132").Value = 3979.2856
$ws.Range("K132").Value = 5981.25
$ws.Range("L132").Value = 11937.8568
$ws.Range("M132").Value = -3451.25
$ws.Range("N132").Value = -16997.8568
$ws.Range("H134").Value = 4740.5
$ws.Range("I134").Value = 3492.5
$ws.Range("J134").Value = 5156.5
$ws.Range("K134").Value = 10477.5
$ws.Range("L134").Value = 15469.5
$ws.Range("M134").Value = -7942.5
$ws.Range("H136").Value = 3469.2942
$ws.Range("I136").Value = 1935.1111
$ws.Range("J136").Value = 5195.25
$ws.Range("K136").Value = 5805.3333
$ws.Range("L136").Value = 15585.75
$ws.Range("M136").Value = -3255.3333

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 1146.7
$ws.Range("I14").Value = 1146.7
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 3440.1
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -3267.1
$ws.Range("H29").Value = 333936.84
$ws.Range("I29").Value = 667294.7
$ws.Range("J29").Value = 579
$ws.Range("K29").Value = 2001884.1
$ws.Range("L29").Value = 1737
$ws.Range("M29").Value = -2001607.1
$ws.Range("N29").Value = -2291
$ws.Range("H107").Value = 1235.4546
$ws.Range("I107").Value = 1156.5
$ws.Range("J107").Value = 1253
$ws.Range("K107").Value = 3469.5
$ws.Range("L107").Value = 3759
$ws.Range("M107").Value = -1549.5
$ws.Range("N107").Value = -7599
$ws.Range("H117").Value = 1404.3334
$ws.Range("I117").Value = 633.3333
$ws.Range("J117").Value = 1789.8334
$ws.Range("K117").Value = 1899.9999
$ws.Range("L117").Value = 5369.5002
$ws.Range("M117").Value = 1542.0001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3007.3333
$ws.Range("I80").Value = 3296.3333
$ws.Range("J80").Value = 2429.3333
$ws.Range("K80").Value = 3296.3333
$ws.Range("L80").Value = 2429.3333
$ws.Range("M80").Value = -2298.3333
$ws.Range("H83").Value = 3007.3333
$ws.Range("I83").Value = 3296.3333
$ws.Range("J83").Value = 2429.3333
$ws.Range("K83").Value = 16481.6665
$ws.Range("L83").Value = 12146.6665
$ws.Range("M83").Value = -11489.6665
$ws.Range("H97").Value = 1873.8
$ws.Range("I97").Value = 1891.1333
$ws.Range("J97").Value = 1821.8
$ws.Range("K97").Value = 1891.1333
$ws.Range("L97").Value = 1821.8
$ws.Range("M97").Value = -1395.1333
$ws.Range("N97").Value = -2813.8
$ws.Range("H122").Value = 692919.4399999999
$ws.Range("I122").Value = 116304.78
$ws.Range("J122").Value = 1434281.1
$ws.Range("K122").Value = 348914.34
$ws.Range("L122").Value = 4302843.300000001
$ws.Range("M122").Value = -346464.34
$ws.Range("N122").Value = -4307743.300000001
$ws.Range("H132").Value = 2990.1936
$ws.Range("I132").Value = 2636.76
$ws.Range("J132").Value = 4462.8335
$ws.Range("K132").Value = 7910.280000000001
$ws.Range("L132").Value = 13388.5005
$ws.Range("M132").Value = -5380.280000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 259.16666
$ws.Range("I93").Value = 211
$ws.Range("J93").Value = 500
$ws.Range("K93").Value = 211
$ws.Range("L93").Value = 500
$ws.Range("M93").Value = 1037
$ws.Range("H136").Value = 3169.238
$ws.Range("I136").Value = 3030.8333
$ws.Range("J136").Value = 3999.6667
$ws.Range("K136").Value = 9092.499899999999
$ws.Range("L136").Value = 11999.0001
$ws.Range("M136").Value = -6542.499899999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H60").Value = 99997
$ws.Range("I60").Value = 99994
$ws.Range("J60").Value = 100000
$ws.Range("K60").Value = 99994
$ws.Range("L60").Value = 100000
$ws.Range("M60").Value = -99172
$ws.Range("N60").Value = -101644
$ws.Range("H81").Value = 1378.375
$ws.Range("I81").Value = 1399.1
$ws.Range("J81").Value = 1343.8334
$ws.Range("K81").Value = 2798.2
$ws.Range("L81").Value = 2687.6668
$ws.Range("M81").Value = -1737.2
$ws.Range("N81").Value = -4809.6668
$ws.Range("H84").Value = 1378.375
$ws.Range("I84").Value = 1399.1
$ws.Range("J84").Value = 1343.8334
$ws.Range("K84").Value = 13991
$ws.Range("L84").Value = 13438.334
$ws.Range("M84").Value = -8687
$ws.Range("N84").Value = -24046.334
$ws.Range("H132").Value = 1367.3889
$ws.Range("I132").Value = 1448.4445
$ws.Range("J132").Value = 1286.3334
$ws.Range("K132").Value = 4345.333500000001
$ws.Range("L132").Value = 3859.0002
$ws.Range("M132").Value = -1815.333500000001
$ws.Range("N132").Value = -8919.0002
$ws.Range("H136").Value = 68865.07000000001
$ws.Range("I136").Value = 1484
$ws.Range("J136").Value = 254163
$ws.Range("K136").Value = 4452
$ws.Range("L136").Value = 762489
$ws.Range("M136").Value = -1902
$ws.Range("N136").Value = -767589

